$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.023.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.455.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.454.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("E11").Value = "  -0.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("E13").Value = "  -7.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.888.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.881.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.370.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("E22").Value = "  +4.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -2.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.381"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.89%  "

$ws.Range("E29").Value = "  +4.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0735"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("E38").Value = "  +5.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.94%  "

$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "256.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0492"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +1.58%  "

$ws.Range("E51").Value = "  +0.35%  "
